$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet): update "想去人数" (interest count) figures
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1306
$ws1.Range("F3").Value = 2821

# Sheet "全部类型" (fourth sheet): same events repeated, update matching figures
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1306
$ws4.Range("F4").Value = 2821
